# PAE - MultiHiper - inclusão da nota do simulado sem nome
#
# The original sheet lists one exam-score row per student (columns:
# A = Nome, B = Identificação de Usuário, C = Simulado). This adds a
# trailing row for the submission that came in without a name attached
# ("Sem Nome"), with no user id, and highlights its passing score using
# the same "Good"/"Bom" look Excel's built-in conditional formatting
# already uses for scores >= 5 elsewhere in the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 63

# New "Bom" cell style (Excel's built-in "Good" style, id 26: green text
# on a light-green fill) so the added score cell is flagged the same way
# passing scores already are via the sheet's conditional formatting.
$bom = $wb.Styles.Add("Bom")
$bom.Font.Color = 24832        # RGB(0, 97, 0)   -> FF006100
$bom.Interior.Color = 13561798 # RGB(198,239,206) -> FFC6EFCE

$ws.Cells.Item($row, 1).Value = "Sem Nome"
$ws.Cells.Item($row, 3).Value = 8.1
$ws.Cells.Item($row, 3).Style = "Bom"
